# DiscountRules.xlsx update
#
# Effective change (decoded from the OOXML diff): a new "CONDITION" column is
# inserted immediately before the old column E, pushing the previous E -> F
# and the previous F -> G. Two new shared strings are introduced ("next vv"
# in the header/body and "03102025 0937" used for a brand-new row 27).
#
# Because Range.Insert() on this host re-derives the whole shared-string
# table (and because plain Value="" assignment clears a cell instead of
# storing a literal empty string), every target cell below is written
# explicitly instead of relying on a structural column insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextCell {
    # Writes a literal text value - including a true empty string, which is
    # distinct from a blank cell - using a quote-prefixed entry ('<text>) so
    # numeric-looking strings ("25", "30", ...) stay text instead of being
    # coerced into numeric cells. The style is reset back to Normal right
    # after so the cell keeps using the sheet's default (s="0") format.
    param($cell, [string]$text)
    $cell.Formula = "'" + $text
    $cell.Style = "Normal"
}

function Set-BlankCell {
    # Materialises a present-but-empty cell (no <v>, no t="s") at the
    # default style, matching the workbook's existing blank placeholder
    # cells (e.g. the source F25/F26 cells).
    param($cell)
    $cell.Formula = "'x"
    $cell.ClearContents()
    $cell.Style = "Normal"
}

# --- Row 18 (header row) ---------------------------------------------------
Set-TextCell $ws.Cells.Item(18, 5) "CONDITION"   # E18 (new column)
Set-TextCell $ws.Cells.Item(18, 6) "ACTION"      # F18 (was E18)
Set-TextCell $ws.Cells.Item(18, 7) "ACTION"      # G18 (was F18)

# --- Row 19 ------------------------------------------------------------
Set-TextCell $ws.Cells.Item(19, 5) "next vv"     # E19 (new column)
Set-TextCell $ws.Cells.Item(19, 6) "jnn"         # F19 (was E19)
Set-TextCell $ws.Cells.Item(19, 7) "next"        # G19 (was F19)

# --- Row 20 ------------------------------------------------------------
Set-TextCell $ws.Cells.Item(20, 5) ""            # E20 (new column)
Set-TextCell $ws.Cells.Item(20, 6) "25"          # F20 (was E20)
Set-TextCell $ws.Cells.Item(20, 7) "25"          # G20 (was F20)

# --- Row 21 ------------------------------------------------------------
Set-TextCell $ws.Cells.Item(21, 5) ""            # E21 (new column)
Set-TextCell $ws.Cells.Item(21, 6) "30"          # F21 (was E21)
Set-TextCell $ws.Cells.Item(21, 7) "30"          # G21 (was F21)

# --- Row 22 ------------------------------------------------------------
Set-TextCell $ws.Cells.Item(22, 5) ""            # E22 (new column)
Set-TextCell $ws.Cells.Item(22, 6) "35"          # F22 (was E22)
Set-TextCell $ws.Cells.Item(22, 7) "35"          # G22 (was F22)

# --- Row 23 ------------------------------------------------------------
Set-TextCell $ws.Cells.Item(23, 5) ""            # E23 (new column)
Set-TextCell $ws.Cells.Item(23, 6) "40"          # F23 (was E23)
Set-TextCell $ws.Cells.Item(23, 7) "40"          # G23 (was F23)

# --- Row 24 ------------------------------------------------------------
Set-TextCell $ws.Cells.Item(24, 5) ""            # E24 (new column)
Set-TextCell $ws.Cells.Item(24, 6) "45"          # F24 (was E24)
Set-TextCell $ws.Cells.Item(24, 7) "45"          # G24 (was F24)

# --- Row 25 ------------------------------------------------------------
Set-TextCell $ws.Cells.Item(25, 5) ""            # E25 (new column; was truly blank)
# F25 keeps its pre-existing blank cell untouched (no value)
Set-BlankCell $ws.Cells.Item(25, 7)              # G25 (new blank cell)

# --- Row 26 ------------------------------------------------------------
Set-TextCell $ws.Cells.Item(26, 5) ""            # E26 (new column; was truly blank)
# F26 keeps its pre-existing blank cell untouched (no value)
Set-BlankCell $ws.Cells.Item(26, 7)              # G26 (new blank cell)

# --- Row 27 (brand-new row) --------------------------------------------
Set-TextCell $ws.Cells.Item(27, 1) "03102025 0937"   # A27
Set-BlankCell $ws.Cells.Item(27, 2)                  # B27
Set-BlankCell $ws.Cells.Item(27, 3)                  # C27
Set-BlankCell $ws.Cells.Item(27, 4)                  # D27
Set-TextCell $ws.Cells.Item(27, 5) "03102025 0937"   # E27
Set-BlankCell $ws.Cells.Item(27, 6)                  # F27
Set-BlankCell $ws.Cells.Item(27, 7)                  # G27
